# Change start/end time input for api service tests, e.g. getSensorDataByDeviceId
$wb = $excel.ActiveWorkbook

# --- getSensorDataBySensorId: update endTime (C2) / startTime (E2) values ---
$wsSensorDataBySensorId = $wb.Worksheets.Item("getSensorDataBySensorId")
$wsSensorDataBySensorId.Range("C2").Value = "2021-6-17 03:00:00"
$wsSensorDataBySensorId.Range("E2").Value = "2021-1-27 04:00:00"

# --- getSensorDataByDeviceId: update endTime (C2) / startTime (F2) values ---
$wsSensorDataByDeviceId = $wb.Worksheets.Item("getSensorDataByDeviceId")
$wsSensorDataByDeviceId.Range("C2").Value = "2021-06-17 03:00:00"
$wsSensorDataByDeviceId.Range("F2").Value = "2021-1-27 04:00:00"

# --- getKpiDataByDeviceId: update endTime (C2) / startTime (F2) values ---
$wsKpiDataByDeviceId = $wb.Worksheets.Item("getKpiDataByDeviceId")
$wsKpiDataByDeviceId.Range("C2").Value = "2021-06-17 15:00:00"
$wsKpiDataByDeviceId.Range("F2").Value = "2020-12-27 04:00:00"

# --- Update the active cell / selection on each touched sheet, and which sheet / tab is active ---
$wsSensorDataByDeviceId.Activate()
$wsSensorDataByDeviceId.Range("F6").Select() | Out-Null

$wsKpiDataByDeviceId.Activate()
$wsKpiDataByDeviceId.Range("E4").Select() | Out-Null

$wsSensorDataBySensorId.Activate()
$wsSensorDataBySensorId.Range("F4").Select() | Out-Null
